$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new product ("KETOLAC 30MG/2ML 5 AMP. FOR I.M./I.V. INF.") was added to
# the report, in its alphabetically-sorted position between "IVYROSPAN SYRUP
# 100 ML" (row 9) and "MAXILASE 200 CEIP UNIT/ML SYRUP 100ML" (old row 10).
# That pushes every following data row down by one, the grand-total row
# (previously 16) down to 17, and the footer row (previously 17) down to 18.
# ---------------------------------------------------------------------------

# Insert a new blank row at row 10, shifting rows 10-17 down to 11-18
# (mergeCells and row content move with it automatically).
$ws.Rows(10).Insert()

# Clone the formatting (styles + row height) of the row directly below
# (the old row 10, now shifted to row 11) into the freshly inserted row 10,
# so the new row matches the look of every other product row.
$ws.Range("A11:N11").Copy()
$ws.Range("A10:N10").PasteSpecial(-4122)
$ws.Rows(10).RowHeight = 25.5

# Recreate the same merged ranges used by every other product row.
$ws.Range("B10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()

# Fill in the new row's data.
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "KETOLAC 30MG/2ML 5 AMP. FOR I.M./I.V. INF."
$ws.Range("H10").Value = "3:2"
$ws.Range("L10").Value = 12
$ws.Range("N10").Value = "0:0"

# The serial numbers (column A) in the shifted-down rows stay at their old
# value after the insert, so bump each one up by one to keep them sequential.
$ws.Range("A11").Value = 8
$ws.Range("A12").Value = 9
$ws.Range("A13").Value = 10
$ws.Range("A14").Value = 11
$ws.Range("A15").Value = 12
$ws.Range("A16").Value = 13

# The "سرنجات 3 سم" row (now row 15 after the shift) also got its own data
# refresh for this period.
$ws.Range("H15").Value = "-1:0"
$ws.Range("L15").Value = 16
$ws.Range("N15").Value = "8:0"

# The grand-total row (now row 17) is updated to reflect the new total.
$ws.Range("K17").Value = 739

# Row heights are keyed off the row position, not the content that moved
# into them, so reassert the exact heights for every row from the inserted
# one down through the (shifted) total/footer rows.
$ws.Rows(10).RowHeight = 25.5
$ws.Rows(11).RowHeight = 24.75
$ws.Rows(12).RowHeight = 25.5
$ws.Rows(13).RowHeight = 25.5
$ws.Rows(14).RowHeight = 24.75
$ws.Rows(15).RowHeight = 25.5
$ws.Rows(16).RowHeight = 24.75
$ws.Rows(17).RowHeight = 26.25
$ws.Rows(18).RowHeight = 16.5
